$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value2 = 10516.47
$ws.Cells.Item(32, 9).Value2 = 7765.758
$ws.Cells.Item(32, 11).Value2 = 7765.758
$ws.Cells.Item(32, 13).Value2 = -7478.758

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value2 = 2217.6829
$ws.Cells.Item(61, 9).Value2 = 2019.6086
$ws.Cells.Item(61, 10).Value2 = 2470.7778
$ws.Cells.Item(61, 11).Value2 = 2019.6086
$ws.Cells.Item(61, 12).Value2 = 2470.7778
$ws.Cells.Item(61, 13).Value2 = -1807.6086
$ws.Cells.Item(61, 14).Value2 = -2894.7778

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(74, 8).Value2 = 8622065
$ws.Cells.Item(74, 9).Value2 = 12821682
$ws.Cells.Item(74, 10).Value2 = 1800.6842
$ws.Cells.Item(74, 11).Value2 = 12821682
$ws.Cells.Item(74, 12).Value2 = 1800.6842
$ws.Cells.Item(74, 13).Value2 = -12820808
$ws.Cells.Item(74, 14).Value2 = -3548.6842

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(77, 8).Value2 = 8622065
$ws.Cells.Item(77, 9).Value2 = 12821682
$ws.Cells.Item(77, 10).Value2 = 1800.6842
$ws.Cells.Item(77, 11).Value2 = 64108410
$ws.Cells.Item(77, 12).Value2 = 9003.421
$ws.Cells.Item(77, 13).Value2 = -64104042
$ws.Cells.Item(77, 14).Value2 = -17739.421

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(136, 8).Value2 = 2217.6829
$ws.Cells.Item(136, 9).Value2 = 2019.6086
$ws.Cells.Item(136, 10).Value2 = 2470.7778
$ws.Cells.Item(136, 11).Value2 = 6058.825800000001
$ws.Cells.Item(136, 12).Value2 = 7412.3334
$ws.Cells.Item(136, 13).Value2 = -3508.825800000001
$ws.Cells.Item(136, 14).Value2 = -12512.3334

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value2 = 9262607
$ws.Cells.Item(31, 9).Value2 = 17858352
$ws.Cells.Item(31, 10).Value2 = 5652.154
$ws.Cells.Item(31, 11).Value2 = 17858352
$ws.Cells.Item(31, 12).Value2 = 5652.154
$ws.Cells.Item(31, 13).Value2 = -17858057
$ws.Cells.Item(31, 14).Value2 = -6242.154

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(34, 8).Value2 = 9262607
$ws.Cells.Item(34, 9).Value2 = 17858352
$ws.Cells.Item(34, 10).Value2 = 5652.154
$ws.Cells.Item(34, 11).Value2 = 17858352
$ws.Cells.Item(34, 12).Value2 = 5652.154
$ws.Cells.Item(34, 13).Value2 = -17858150
$ws.Cells.Item(34, 14).Value2 = -6056.154

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(3, 8).Value2 = 4564.913
$ws.Cells.Item(3, 9).Value2 = 3745
$ws.Cells.Item(3, 10).Value2 = 5002.2
$ws.Cells.Item(3, 11).Value2 = 11235
$ws.Cells.Item(3, 12).Value2 = 15006.6
$ws.Cells.Item(3, 13).Value2 = -11123
$ws.Cells.Item(3, 14).Value2 = -15230.6

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 8).Value2 = 397.04544
$ws.Cells.Item(5, 9).Value2 = 249.21053
$ws.Cells.Item(5, 10).Value2 = 1333.3334
$ws.Cells.Item(5, 11).Value2 = 747.63159
$ws.Cells.Item(5, 12).Value2 = 4000.0002
$ws.Cells.Item(5, 13).Value2 = -635.63159
$ws.Cells.Item(5, 14).Value2 = -4224.0002

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(8, 8).Value2 = 171.92857
$ws.Cells.Item(8, 9).Value2 = 171.92857
$ws.Cells.Item(8, 11).Value2 = 515.78571
$ws.Cells.Item(8, 13).Value2 = -376.78571

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(75, 8).Value2 = 801.8333
$ws.Cells.Item(75, 9).Value2 = 637.6667
$ws.Cells.Item(75, 10).Value2 = 966
$ws.Cells.Item(75, 11).Value2 = 1913.0001
$ws.Cells.Item(75, 12).Value2 = 2898
$ws.Cells.Item(75, 13).Value2 = -915.0001
$ws.Cells.Item(75, 14).Value2 = -4894

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(78, 8).Value2 = 801.8333
$ws.Cells.Item(78, 9).Value2 = 637.6667
$ws.Cells.Item(78, 10).Value2 = 966
$ws.Cells.Item(78, 11).Value2 = 5739.0003
$ws.Cells.Item(78, 12).Value2 = 8694
$ws.Cells.Item(78, 13).Value2 = -747.0002999999997
$ws.Cells.Item(78, 14).Value2 = -18678

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(103, 8).Value2 = 1878.6
$ws.Cells.Item(103, 9).Value2 = 394.125
$ws.Cells.Item(103, 11).Value2 = 1182.375
$ws.Cells.Item(103, 13).Value2 = -303.375

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(114, 8).Value2 = 902.65216
$ws.Cells.Item(114, 9).Value2 = 342.2857
$ws.Cells.Item(114, 10).Value2 = 1147.8125
$ws.Cells.Item(114, 11).Value2 = 1026.8571
$ws.Cells.Item(114, 12).Value2 = 3443.4375
$ws.Cells.Item(114, 13).Value2 = 2227.1429
$ws.Cells.Item(114, 14).Value2 = -9951.4375

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(117, 8).Value2 = 2764.2666
$ws.Cells.Item(117, 10).Value2 = 4118.222
$ws.Cells.Item(117, 12).Value2 = 12354.666
$ws.Cells.Item(117, 14).Value2 = -19238.666

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(129, 8).Value2 = 674.2857
$ws.Cells.Item(129, 9).Value2 = 674.2857
$ws.Cells.Item(129, 10).Value2 = 0
$ws.Cells.Item(129, 11).Value2 = 2022.8571
$ws.Cells.Item(129, 12).Value2 = 0
$ws.Cells.Item(129, 13).Value2 = 2977.1429
$ws.Cells.Item(129, 14).ClearContents()

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(133, 8).Value2 = 7042.5713
$ws.Cells.Item(133, 9).Value2 = 2000
$ws.Cells.Item(133, 10).Value2 = 7883
$ws.Cells.Item(133, 11).Value2 = 6000
$ws.Cells.Item(133, 12).Value2 = 23649
$ws.Cells.Item(133, 13).Value2 = -940
$ws.Cells.Item(133, 14).Value2 = -33769

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(134, 8).Value2 = 4108
$ws.Cells.Item(134, 9).Value2 = 1617.5
$ws.Cells.Item(134, 10).Value2 = 4974.2607
$ws.Cells.Item(134, 11).Value2 = 4852.5
$ws.Cells.Item(134, 12).Value2 = 14922.7821
$ws.Cells.Item(134, 13).Value2 = 217.5
$ws.Cells.Item(134, 14).Value2 = -25062.7821

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(135, 8).Value2 = 397.04544
$ws.Cells.Item(135, 9).Value2 = 249.21053
$ws.Cells.Item(135, 10).Value2 = 1333.3334
$ws.Cells.Item(135, 11).Value2 = 2242.89477
$ws.Cells.Item(135, 12).Value2 = 12000.0006
$ws.Cells.Item(135, 13).Value2 = 292.1052300000001
$ws.Cells.Item(135, 14).Value2 = -17070.0006

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(136, 8).Value2 = 1630.5
$ws.Cells.Item(136, 10).Value2 = 1597.875
$ws.Cells.Item(136, 12).Value2 = 4793.625
$ws.Cells.Item(136, 14).Value2 = -14993.625

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(138, 8).Value2 = 1569.5834
$ws.Cells.Item(138, 9).Value2 = 1319.2858
$ws.Cells.Item(138, 10).Value2 = 1920
$ws.Cells.Item(138, 11).Value2 = 3957.8574
$ws.Cells.Item(138, 12).Value2 = 5760
$ws.Cells.Item(138, 13).Value2 = 1182.1426
$ws.Cells.Item(138, 14).Value2 = -16040

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(139, 8).Value2 = 3185.8276
$ws.Cells.Item(139, 9).Value2 = 1888.2778
$ws.Cells.Item(139, 10).Value2 = 5309.091
$ws.Cells.Item(139, 11).Value2 = 5664.8334
$ws.Cells.Item(139, 12).Value2 = 15927.273
$ws.Cells.Item(139, 13).Value2 = -524.8334000000004
$ws.Cells.Item(139, 14).Value2 = -26207.273

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(113, 8).Value2 = 6668747.5
$ws.Cells.Item(113, 9).Value2 = 12502100
$ws.Cells.Item(113, 10).Value2 = 2059
$ws.Cells.Item(113, 11).Value2 = 12502100
$ws.Cells.Item(113, 12).Value2 = 2059
$ws.Cells.Item(113, 13).Value2 = -12499930
$ws.Cells.Item(113, 14).Value2 = -6399

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(126, 8).Value2 = 12501971
$ws.Cells.Item(126, 9).Value2 = 14707831
$ws.Cells.Item(126, 10).Value2 = 2100
$ws.Cells.Item(126, 11).Value2 = 44123493
$ws.Cells.Item(126, 12).Value2 = 6300
$ws.Cells.Item(126, 13).Value2 = -44121023
$ws.Cells.Item(126, 14).Value2 = -11240

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value2 = 2504.282
$ws.Cells.Item(40, 9).Value2 = 1965.4615
$ws.Cells.Item(40, 11).Value2 = 1965.4615
$ws.Cells.Item(40, 13).Value2 = -1829.4615

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(136, 8).Value2 = 3060.724
$ws.Cells.Item(136, 9).Value2 = 2800.0625
$ws.Cells.Item(136, 10).Value2 = 3381.5386
$ws.Cells.Item(136, 11).Value2 = 8400.1875
$ws.Cells.Item(136, 12).Value2 = 10144.6158
$ws.Cells.Item(136, 13).Value2 = -5850.1875
$ws.Cells.Item(136, 14).Value2 = -15244.6158

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(108, 8).Value2 = 34375.145
$ws.Cells.Item(108, 10).Value2 = 34375.145
$ws.Cells.Item(108, 12).Value2 = 34375.145
$ws.Cells.Item(108, 14).Value2 = -42055.145

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(136, 8).Value2 = 4873.4165
$ws.Cells.Item(136, 9).Value2 = 1564.4348
$ws.Cells.Item(136, 10).Value2 = 10727.77
$ws.Cells.Item(136, 11).Value2 = 4693.3044
$ws.Cells.Item(136, 12).Value2 = 32183.31
$ws.Cells.Item(136, 13).Value2 = -2143.3044
$ws.Cells.Item(136, 14).Value2 = -37283.31
